# Edit: split the "Participation increased by 15% ... " sentence on slide 15
# so that "15%" becomes "14%", matching the author's commit which re-typed
# that portion of the sentence (and so it now lives in its own run).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$target = "by 15% "
$idx = $full.IndexOf($target)
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, $target.Length)
    $chars.Text = "by 14% "
}
